$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C width (11.7109375 -> 12.7109375 character-width units).
# Excel's ColumnWidth property is quantized to whole pixels (stored width =
# (round(ColumnWidth*6)+5)/6), so we pick the ColumnWidth that rounds to the
# closest achievable stored width (12.666666666666666, pixel=71).
$ws.Columns.Item(3).ColumnWidth = 11.833333333333334

# Update cell values
$ws.Range("A1").Value = 160.95569524474391
$ws.Range("B1").Value = 6.4697087742661701
$ws.Range("C1").Value = 0.70978473581213297
